$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 141; existing rows 141-169 shift down to 142-170
# (Excel copies the formatting of the row being pushed down, which keeps
# column D's date number format intact for the new row.)
$ws.Rows.Item(141).Insert()

# Populate the newly inserted row 141 with the new weekly data point.
$ws.Cells.Item(141, 1).Value = 7
$ws.Cells.Item(141, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(141, 3).Value = "Ñuble"
$ws.Cells.Item(141, 4).Value = [datetime]"2021-10-07"
$ws.Cells.Item(141, 5).Value = 16
$ws.Cells.Item(141, 6).Value = 100112023
$ws.Cells.Item(141, 7).Value = "Brócoli"
$ws.Cells.Item(141, 8).Value = "Sin especificar"
$ws.Cells.Item(141, 9).Value = "Primera"
$ws.Cells.Item(141, 10).Value = 300
$ws.Cells.Item(141, 11).Value = 700
$ws.Cells.Item(141, 12).Value = 750
$ws.Cells.Item(141, 13).Value = 725
$ws.Cells.Item(141, 14).Value = "$/unidad"
$ws.Cells.Item(141, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(141, 16).Value = 725
$ws.Cells.Item(141, 17).Value = 1
$ws.Cells.Item(141, 18).Value = "Hortaliza"
